# Update data in spreadsheet and Methodology document
# Applies the data-driven changes to the "Energy" worksheet of the workbook:
#  - switch several duration inputs (column D) from literal decimal values
#    to equivalent fraction formulas (e.g. 0.16 -> =1/6), which also nudges
#    their dependent "Energy used (Wh)" results (column E) to the
#    higher-precision figures
#  - consolidate the per-row "=B*D" formulas in rows 35-45 into one shared
#    formula group
#  - refresh view/selection state to match where the author was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Energy")

# --- Row 19 (Microwave): D19 formula changes from 1/60*5 to the equivalent 1/12 ---
$ws.Range("D19").Formula = "=1/12"

# --- Row 28 (Vaccum cleaner (hoover)): D28 literal 0.16 -> formula 1/6 ---
$ws.Range("D28").Formula = "=1/6"
$ws.Range("D28").NumberFormat = "0.00"

# --- Row 33 (Clothes iron): D33 literal 0.16 -> formula 1/6 ---
$ws.Range("D33").Formula = "=1/6"
$ws.Range("D33").NumberFormat = "0.00"

# --- Row 35 (Hairdryer): D35 literal 0.08 -> formula 1/12 ---
$ws.Range("D35").Formula = "=1/12"
$ws.Range("D35").NumberFormat = "0.00"

# E35:E45 become one shared formula group (B*D per row), replacing the
# individual (but equivalent) per-row formulas that used to live there.
$ws.Range("E35:E45").Formula = "=B35*D35"
$ws.Range("E35").NumberFormat = "0.00"
$ws.Range("E36").NumberFormat = "0"

# --- Row 36 (Electric shower): D36 literal 0.16 -> formula 1/6 ---
$ws.Range("D36").Formula = "=1/6"
$ws.Range("D36").NumberFormat = "0.00"

# --- Row 37 (Gas-powered shower): D37 literal 0.16 -> formula 1/6 ---
$ws.Range("D37").Formula = "=1/6"
$ws.Range("D37").NumberFormat = "0.00"

# --- Column widths: D becomes a new custom-width column, E widens slightly ---
$ws.Columns.Item(4).ColumnWidth = 10.92
$ws.Columns.Item(5).ColumnWidth = 12.75

# --- View / selection state: scroll down a little and land on D23 ---
$ws.Range("D23").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
